$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 날짜 (date) column A: 2025-12-05 -> 2025-12-06 for all data rows (2-5).
# Leading apostrophe forces text storage so the value stays a shared string
# (matching the source data) instead of being auto-parsed into a date serial.
$ws.Range("A2").Value = "'2025-12-06"
$ws.Range("A3").Value = "'2025-12-06"
$ws.Range("A4").Value = "'2025-12-06"
$ws.Range("A5").Value = "'2025-12-06"

# 최종점수 (final score) column K updates.
$ws.Range("K2").Value = 59.5
$ws.Range("K3").Value = 55.5
$ws.Range("K4").Value = 53.5
$ws.Range("K5").Value = 52.5

# MACRO_SCORE column N updates (same new score for all rows).
$ws.Range("N2").Value = 51.54219175917372
$ws.Range("N3").Value = 51.54219175917372
$ws.Range("N4").Value = 51.54219175917372
$ws.Range("N5").Value = 51.54219175917372
